# Daily attendance processing - 2026-01-23 10:10:15
#
# Normalises the "Recorded By" column (G): whichever token(s) in the
# comma-separated list equal "system" (case-insensitive) are moved to the
# front, preserving the relative order of the remaining tokens. Rows with
# no "system" token are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text.Length -eq 0) { continue }
    if ($text.IndexOf(',') -lt 0) { continue }

    $parts = $text -split ','
    $systemParts = @()
    $restParts = @()
    foreach ($p in $parts) {
        $trimmedPart = $p.Trim()
        if ($trimmedPart.ToLower() -eq 'system') {
            $systemParts += $trimmedPart
        } else {
            $restParts += $trimmedPart
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    $newParts = $systemParts + $restParts
    $newVal = [string]::Join(', ', $newParts)

    if ($newVal -ne $text) {
        $cell.Value = $newVal
    }
}
